$wb = $excel.ActiveWorkbook

# --- January 2014 sheet: add the three new timesheet rows ---
$wsJan = $wb.Worksheets.Item("January 2014")

$wsJan.Range("A14").Value = (Get-Date -Year 2014 -Month 1 -Day 17 -Hour 0 -Minute 0 -Second 0)
$wsJan.Range("B14").Value = "Resource View & Listings"
$wsJan.Range("C14").Value = 4.5
$wsJan.Range("D14").Value = "Riaan Bekker"

$wsJan.Range("A15").Value = (Get-Date -Year 2014 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0)
$wsJan.Range("B15").Value = "ResourceView, Listings, Uploads & Search"
$wsJan.Range("C15").Value = 10
$wsJan.Range("D15").Value = "Riaan Bekker"

$wsJan.Range("A16").Value = (Get-Date -Year 2014 -Month 1 -Day 19 -Hour 0 -Minute 0 -Second 0)
$wsJan.Range("B16").Value = "ResourceView, Listings, Uploads & Search"
$wsJan.Range("C16").Value = 9
$wsJan.Range("D16").Value = "Riaan Bekker"

# --- Update the selection (active cell) on every sheet ---
$wsNov = $wb.Worksheets.Item("November 2013")
$wsNov.Range("C20").Select()

$wsDec = $wb.Worksheets.Item("December 2013")
$wsDec.Range("C20").Select()

$wsJan.Range("C17").Select()
